$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 77, shifting the old separator/summary rows down by one.
$ws.Rows(77).Insert()

# The previously-last working entry (row 76) actually ended later than
# recorded - fix its end time.
$ws.Range("E76").Value = 0.52083333333333337

# Fill in the newly inserted row with the missing time entry.
$ws.Range("A77").Value = 2014
$ws.Range("B77").Value = 3
$ws.Range("C77").Value = 18
$ws.Range("D77").Value = 0.54166666666666663
$ws.Range("E77").Value = 0.625
$ws.Range("F77").Formula = "=(E77-D77)*24*60"
$ws.Range("G77").Formula = "=F77/60"

# Keep the selection where the user left off editing.
$ws.Range("A78").Select() | Out-Null
